$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 481.6
$ws.Range("I39").Value = 204
$ws.Range("J39").Value = 666.6667
$ws.Range("K39").Value = 612
$ws.Range("L39").Value = 2000.0001
$ws.Range("M39").Value = -316
$ws.Range("N39").Value = -2592.0001

$ws.Range("H112").Value = 1730.7858
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1730.7858
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 5192.357400000001
$ws.Range("N112").Value = -7408.357400000001
$ws.Range("M112").ClearContents()

$ws.Range("H113").Value = 3082.6667
$ws.Range("I113").Value = 3082.6667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3082.6667
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 171.3332999999998
$ws.Range("N113").ClearContents()

$ws.Range("H124").Value = 62999
$ws.Range("J124").Value = 62999
$ws.Range("L124").Value = 62999
$ws.Range("N124").Value = -72819

$ws.Range("H129").Value = 909.5862
$ws.Range("J129").Value = 852.75
$ws.Range("L129").Value = 2558.25
$ws.Range("N129").Value = -12558.25

$ws.Range("H138").Value = 1665.8125
$ws.Range("I138").Value = 1418.9259
$ws.Range("J138").Value = 2999
$ws.Range("K138").Value = 4256.7777
$ws.Range("L138").Value = 8997
$ws.Range("M138").Value = 883.2223000000004
$ws.Range("N138").Value = -19277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 34484316
$ws.Range("I45").Value = 45455850
$ws.Range("J45").Value = 2352
$ws.Range("K45").Value = 45455850
$ws.Range("L45").Value = 2352
$ws.Range("M45").Value = -45455473
$ws.Range("N45").Value = -3106

$ws.Range("H122").Value = 1651.2273
$ws.Range("I122").Value = 1517.2106
$ws.Range("K122").Value = 4551.6318
$ws.Range("M122").Value = -2101.6318

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 234
$ws.Range("I5").Value = 97.59999999999999
$ws.Range("J5").Value = 575
$ws.Range("K5").Value = 97.59999999999999
$ws.Range("L5").Value = 575
$ws.Range("M5").Value = 15.40000000000001
$ws.Range("N5").Value = -801

$ws.Range("H7").Value = 820
$ws.Range("I7").Value = 900
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 500
$ws.Range("M7").Value = -787
$ws.Range("N7").Value = -726

$ws.Range("H135").Value = 69632.5
$ws.Range("J135").Value = 69632.5
$ws.Range("L135").Value = 69632.5
$ws.Range("N135").Value = -79772.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 15000500
$ws.Range("I2").Value = 30000000
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 30000000
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -29999887
$ws.Range("N2").Value = -1226

$ws.Range("H11").Value = 30000
$ws.Range("J11").Value = 30000
$ws.Range("L11").Value = 30000
$ws.Range("N11").Value = -30280

$ws.Range("H31").Value = 3404442
$ws.Range("I31").Value = 1273.3024
$ws.Range("J31").Value = 6065101
$ws.Range("K31").Value = 1273.3024
$ws.Range("L31").Value = 6065101
$ws.Range("M31").Value = -978.3024
$ws.Range("N31").Value = -6065691

$ws.Range("H34").Value = 3404442
$ws.Range("I34").Value = 1273.3024
$ws.Range("J34").Value = 6065101
$ws.Range("K34").Value = 1273.3024
$ws.Range("L34").Value = 6065101
$ws.Range("M34").Value = -1071.3024
$ws.Range("N34").Value = -6065505

$ws.Range("H45").Value = 8955.666999999999
$ws.Range("I45").Value = 6867
$ws.Range("K45").Value = 6867
$ws.Range("M45").Value = -6274

$ws.Range("H122").Value = 80743
$ws.Range("I122").Value = 100695.414
$ws.Range("J122").Value = 933.3333
$ws.Range("K122").Value = 302086.242
$ws.Range("L122").Value = 2799.9999
$ws.Range("M122").Value = -299636.242
$ws.Range("N122").Value = -7699.9999

$ws.Range("H132").Value = 132647.36
$ws.Range("I132").Value = 2924.5
$ws.Range("J132").Value = 206774.72
$ws.Range("K132").Value = 8773.5
$ws.Range("L132").Value = 620324.16
$ws.Range("M132").Value = -6243.5
$ws.Range("N132").Value = -625384.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1300.8206
$ws.Range("I68").Value = 1141
$ws.Range("J68").Value = 1348.7667
$ws.Range("K68").Value = 3423
$ws.Range("L68").Value = 4046.300099999999
$ws.Range("M68").Value = -2612
$ws.Range("N68").Value = -5668.300099999999

$ws.Range("H71").Value = 1300.8206
$ws.Range("I71").Value = 1141
$ws.Range("J71").Value = 1348.7667
$ws.Range("K71").Value = 10269
$ws.Range("L71").Value = 12138.9003
$ws.Range("M71").Value = -6213
$ws.Range("N71").Value = -20250.9003

$ws.Range("H113").Value = 3946.375
$ws.Range("I113").Value = 5442.524
$ws.Range("J113").Value = 1090.091
$ws.Range("K113").Value = 16327.572
$ws.Range("L113").Value = 3270.273
$ws.Range("M113").Value = -14157.572
$ws.Range("N113").Value = -7610.272999999999

$ws.Range("H121").Value = 207304.06
$ws.Range("I121").Value = 99.5
$ws.Range("J121").Value = 239181.69
$ws.Range("K121").Value = 298.5
$ws.Range("L121").Value = 717545.0700000001
$ws.Range("M121").Value = 1011.5
$ws.Range("N121").Value = -720165.0700000001

$ws.Range("H131").Value = 3653.3408
$ws.Range("I131").Value = 14775.571
$ws.Range("J131").Value = 1549.1351
$ws.Range("K131").Value = 44326.713
$ws.Range("L131").Value = 4647.4053
$ws.Range("M131").Value = -39286.713
$ws.Range("N131").Value = -14727.4053

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 854.9545000000001
$ws.Range("I102").Value = 763.3570999999999
$ws.Range("J102").Value = 1015.25
$ws.Range("K102").Value = 763.3570999999999
$ws.Range("L102").Value = 1015.25
$ws.Range("M102").Value = 858.6429000000001
$ws.Range("N102").Value = -4259.25

$ws.Range("H122").Value = 1449.2778
$ws.Range("I122").Value = 1546.6923
$ws.Range("J122").Value = 1196
$ws.Range("K122").Value = 4640.0769
$ws.Range("L122").Value = 3588
$ws.Range("M122").Value = -2190.0769
$ws.Range("N122").Value = -8488

$ws.Range("H126").Value = 30334.715
$ws.Range("I126").Value = 133670.67
$ws.Range("J126").Value = 2152.182
$ws.Range("K126").Value = 401012.01
$ws.Range("L126").Value = 6456.545999999999
$ws.Range("M126").Value = -398542.01
$ws.Range("N126").Value = -11396.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 71432440
$ws.Range("J7").Value = 5801.4287
$ws.Range("L7").Value = 5801.4287
$ws.Range("N7").Value = -6025.4287

$ws.Range("H30").Value = 108.666664
$ws.Range("I30").Value = 113
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = 113
$ws.Range("L30").Value = 100
$ws.Range("M30").Value = -5
$ws.Range("N30").Value = -316

$ws.Range("H63").Value = 18000
$ws.Range("J63").Value = 18000
$ws.Range("L63").Value = 18000
$ws.Range("N63").Value = -19498

$ws.Range("H66").Value = 18000
$ws.Range("J66").Value = 18000
$ws.Range("L66").Value = 54000
$ws.Range("N66").Value = -61488

$ws.Range("H122").Value = 2476
$ws.Range("I122").Value = 2432.3333
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 7296.999899999999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -4846.999899999999
$ws.Range("N122").Value = -13900

$ws.Range("H126").Value = 71432440
$ws.Range("J126").Value = 5801.4287
$ws.Range("L126").Value = 17404.2861
$ws.Range("N126").Value = -22344.2861

$ws.Range("H134").Value = 64709.5
$ws.Range("J134").Value = 64709.5
$ws.Range("L134").Value = 64709.5
$ws.Range("N134").Value = -74849.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1839455.1
$ws.Range("I126").Value = 2452265.2
$ws.Range("J126").Value = 1025
$ws.Range("K126").Value = 7356795.600000001
$ws.Range("L126").Value = 3075
$ws.Range("M126").Value = -7354325.600000001
$ws.Range("N126").Value = -8015
